$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ireland")

$ws.Range("C30").Value = 1819
$ws.Range("F30").Value = 19
